# Update Name of Algo
# Applies updated imputed values produced by the RandomForest algorithm
# for terrestrial_mammals / combination_2_ABCDE / ABCE / 20 / seed2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.347000000000003
$ws.Range("E2").Value = 16.6774
$ws.Range("A3").Value = -21.84079999999999
$ws.Range("C3").Value = -10.9668
$ws.Range("E6").Value = 16.61930000000001
$ws.Range("C12").Value = -11.72529999999999
$ws.Range("A14").Value = -21.7476
$ws.Range("A16").Value = -22.08050000000001
$ws.Range("B18").Value = 6.817499999999995
$ws.Range("E19").Value = 16.286
$ws.Range("A21").Value = -20.29989999999997
$ws.Range("A23").Value = -20.79259999999998
$ws.Range("B24").Value = 6.428399999999998
$ws.Range("C24").Value = -11.608
$ws.Range("E24").Value = 16.4892
$ws.Range("A25").Value = -22.07579999999999
$ws.Range("B25").Value = 5.106799999999995
$ws.Range("C25").Value = -13.36709999999999
$ws.Range("A26").Value = -21.05649999999996
$ws.Range("B27").Value = 5.998900000000003
$ws.Range("E27").Value = 16.84959999999999
$ws.Range("A29").Value = -20.87469999999998
$ws.Range("B30").Value = 5.7767
$ws.Range("E30").Value = 15.8403
$ws.Range("B31").Value = 5.617000000000003
$ws.Range("E31").Value = 16.3844
$ws.Range("E33").Value = 16.84480000000001
$ws.Range("B39").Value = 9.603300000000001
$ws.Range("A40").Value = -19.12709999999999
$ws.Range("C41").Value = -12.6862
$ws.Range("B42").Value = 10.4059
$ws.Range("E42").Value = 16.29600000000001
$ws.Range("B48").Value = 4.965400000000004
$ws.Range("C50").Value = -13.10239999999999
$ws.Range("B51").Value = 5.1435
$ws.Range("B52").Value = 5.368499999999998
$ws.Range("A53").Value = -21.41150000000001
$ws.Range("C53").Value = -10.2674
$ws.Range("B55").Value = 6.62909999999999
$ws.Range("E55").Value = 16.64200000000001
$ws.Range("B56").Value = 5.913499999999999
$ws.Range("C56").Value = -12.19149999999999
$ws.Range("A57").Value = -21.99649999999999
$ws.Range("B57").Value = 5.129899999999999
$ws.Range("C57").Value = -12.82169999999999
$ws.Range("C58").Value = -13.96059999999999
$ws.Range("E58").Value = 15.95260000000001
$ws.Range("A59").Value = -22.47480000000001
$ws.Range("B60").Value = 5.5982
$ws.Range("C61").Value = -12.69820000000001
$ws.Range("C63").Value = -11.69599999999999
$ws.Range("C64").Value = -11.73829999999999
$ws.Range("A65").Value = -21.65689999999996
$ws.Range("E65").Value = 17.03690000000001
$ws.Range("A69").Value = -21.55630000000001
$ws.Range("C70").Value = -12.50289999999999
$ws.Range("E70").Value = 16.90240000000001
$ws.Range("C72").Value = -11.7431
$ws.Range("B73").Value = 9.022299999999996
$ws.Range("B74").Value = 9.063799999999997
$ws.Range("E74").Value = 16.6208
$ws.Range("E75").Value = 16.55340000000001
$ws.Range("A79").Value = -20.5501
$ws.Range("A83").Value = -21.90039999999999
$ws.Range("E83").Value = 16.61810000000002
$ws.Range("E84").Value = 16.475
$ws.Range("C86").Value = -13.2579
$ws.Range("E86").Value = 16.55410000000001
$ws.Range("B89").Value = 5.376299999999997
$ws.Range("C89").Value = -10.29570000000001
$ws.Range("B90").Value = 5.696000000000001
$ws.Range("A91").Value = -21.41079999999999
$ws.Range("B92").Value = 5.047099999999991
$ws.Range("A93").Value = -21.15619999999998
$ws.Range("E96").Value = 16.5697
$ws.Range("E97").Value = 16.50520000000001
$ws.Range("C98").Value = -11.6481
$ws.Range("A100").Value = -22.0107
$ws.Range("C100").Value = -13.42019999999999
$ws.Range("C102").Value = -12.493
